$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the bonus placeholder cell to use the custom expression notation
# "[[ ... ]]" instead of the default "${ ... }" notation.
$ws.Range("D4").Value = "[[employee.bonus]]"
